$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab (sheet name changes from "Through 2022-03-01" to "Through 2022-03-02")
$ws.Name = "Through 2022-03-02"

# Update the "March (through 03-01)" label to "March (through 03-02)"
$ws.Range("A4").Value = "March (through 03-02)"

# Update March row (row 4) values for 2018-2022 (columns E-I)
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 6
$ws.Range("I4").Value = 9

# Update Total row (row 5) values for 2018-2022 (columns E-I)
$ws.Range("E5").Value = 139
$ws.Range("F5").Value = 82
$ws.Range("G5").Value = 144
$ws.Range("H5").Value = 348
$ws.Range("I5").Value = 309
